$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newRow = 92

# New daily gold-price entry appended to the table.
$dateValue = "17-12-2025"
$priceValue = "The price of gold in India today is ₹13,451 per gram for 24 karat gold, ₹12,330 per gram for 22 karat gold and ₹10,088 per gram for 18 karat gold (also called 999 gold)."

# Match the formatting used by the previous (last) data row.
$ws.Range("A91:B91").Copy() | Out-Null
$ws.Range("A92:B92").PasteSpecial(-4122) | Out-Null

$ws.Cells.Item($newRow, 1).Value = $dateValue
$ws.Cells.Item($newRow, 2).Value = $priceValue
